$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new TPM run adds a third target cluster ("ECs") for the Rln3->Rxfp2
# pair, placed as the new first data row, and recomputes the specificity
# scores (columns O, P, S, T) for the other two rows now that there are
# three target clusters instead of two. Shift the existing two data rows
# down one position (row2->row3, row3->row4) by copy/paste (this keeps
# the default/no style on the rows, same as before), then fill in the
# new row2 with the ECs record and patch up the recomputed columns.

$ws.Range("A3:T3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial() | Out-Null

$ws.Range("A2:T2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial() | Out-Null

$ws.Application.CutCopyMode = $false

# New row 2: FAPs -> Rln3/Rxfp2 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rln3"
$ws.Range("C2").Value = "Rxfp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.114225
$ws.Range("H2").Value = 0.342675
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.0002143333333333333
$ws.Range("N2").Value = 0.000643
$ws.Range("O2").Value = 0.0005591790590486129
$ws.Range("P2").Value = 0.0005591790590486129
$ws.Range("Q2").Value = 0.000024482225
$ws.Range("R2").Value = 0.000220340025
$ws.Range("S2").Value = 0.0005591790590486129
$ws.Range("T2").Value = 0.0005591790590486129

# Row 3 (previously row 2, target FAPs->FAPs): target cluster text is
# unchanged (still FAPs) but its specificity-derived columns are
# recomputed now that 3 rows exist.
$ws.Range("D3").Value = "FAPs"
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.380202
$ws.Range("N3").Value = 1.140606
$ws.Range("O3").Value = 0.991917558048526
$ws.Range("P3").Value = 0.991917558048526
$ws.Range("Q3").Value = 0.04342857345
$ws.Range("R3").Value = 0.39085716105
$ws.Range("S3").Value = 0.991917558048526
$ws.Range("T3").Value = 0.991917558048526

# Row 4 (previously row 3, target MuSCs): values stay the same except
# the recomputed specificity columns O, P, S, T.
$ws.Range("O4").Value = 0.007523262892425429
$ws.Range("P4").Value = 0.00752326289242543
$ws.Range("S4").Value = 0.007523262892425429
$ws.Range("T4").Value = 0.00752326289242543
